$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7850438356399536
$ws.Range("B1").Value = 2.729622364044189
$ws.Range("C1").Value = 7.698080539703369
$ws.Range("D1").Value = 2.270837306976318
$ws.Range("E1").Value = 1.495089292526245
